$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new value would otherwise be parsed as a number
$protected = $ws.Range("D5,D6,D13,D14,D20,D21,D24,D26,D27,D39,D40,D41,D42,D46,D50,D51")
$protected.NumberFormat = "@"

# Write updated values (text-safe cells)
$ws.Range("D2").Value = "62.072.85"
$ws.Range("E2").Value = "  -0.59%  "
$ws.Range("D3").Value = "2.444.25"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("E5").Value = "  +1.76%  "
$ws.Range("E6").Value = "  -0.74%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  +0.33%  "
$ws.Range("D9").Value = "2.438.95"
$ws.Range("E9").Value = "  +0.33%  "
$ws.Range("E10").Value = "  +1.63%  "
$ws.Range("E11").Value = "  +2.48%  "
$ws.Range("E12").Value = "  -0.18%  "
$ws.Range("E13").Value = "  -2.64%  "
$ws.Range("E14").Value = "  -0.69%  "
$ws.Range("E15").Value = "  +0.58%  "
$ws.Range("D16").Value = "2.884.89"
$ws.Range("E16").Value = "  +0.46%  "
$ws.Range("D17").Value = "62.122.42"
$ws.Range("E17").Value = "  -0.28%  "
$ws.Range("D18").Value = "2.432.67"
$ws.Range("E18").Value = "  -0.26%  "
$ws.Range("E19").Value = "  -3.04%  "
$ws.Range("E20").Value = "  +1.37%  "
$ws.Range("E22").Value = "  -1.04%  "
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("E24").Value = "  -5.58%  "
$ws.Range("E25").Value = "  +0.32%  "
$ws.Range("E26").Value = "  +1.02%  "
$ws.Range("E27").Value = "  -4.33%  "
$ws.Range("D28").Value = "0.0₃0962"
$ws.Range("E28").Value = "  +0.34%  "
$ws.Range("D29").Value = "2.563.51"
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("E30").Value = "  -0.28%  "
$ws.Range("E31").Value = "  -1.31%  "
$ws.Range("E32").Value = "  -0.82%  "
$ws.Range("E33").Value = "  -0.38%  "
$ws.Range("E34").Value = "  +0.17%  "
$ws.Range("E35").Value = "  -2.21%  "
$ws.Range("E36").Value = "  +0.21%  "
$ws.Range("E37").Value = "  -1.45%  "
$ws.Range("E38").Value = "  +0.23%  "
$ws.Range("E39").Value = "  +4.05%  "
$ws.Range("E40").Value = "  -0.35%  "
$ws.Range("E41").Value = "  +0.79%  "
$ws.Range("E42").Value = "  +2.10%  "
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("E45").Value = "  +0.78%  "
$ws.Range("E47").Value = "  -2.50%  "
$ws.Range("D48").Value = "0.0₆0266"
$ws.Range("E48").Value = "  +19.15%  "
$ws.Range("E49").Value = "  +0.63%  "
$ws.Range("E50").Value = "  -0.72%  "
$ws.Range("E51").Value = "  +0.26%  "

# Write updated values (cells protected above so they stay text)
$ws.Range("D5").Value = "583.05"
$ws.Range("D6").Value = "142.51"
$ws.Range("D13").Value = "0.342"
$ws.Range("D14").Value = "26.32"
$ws.Range("D20").Value = "7.21"
$ws.Range("D21").Value = "325.78"
$ws.Range("D24").Value = "1.91"
$ws.Range("D26").Value = "9.12"
$ws.Range("D27").Value = "599.99"
$ws.Range("D39").Value = "152.85"
$ws.Range("D40").Value = "18.38"
$ws.Range("D41").Value = "5.27"
$ws.Range("D42").Value = "43.15"
$ws.Range("D46").Value = "141.79"
$ws.Range("D50").Value = "0.0520"
$ws.Range("D51").Value = "19.76"

# Remove the temporary text formatting so styling matches the original (unstyled) cells
$protected.ClearFormats()
